$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = "Indore"
$ws.Range("C4").Value = "bhopal"
$ws.Range("D4").Value = "07-Mar-2021"
